$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: new "Content Upload" test case
$ws.Range("A21").Value = "Upload_01"
$ws.Range("B21").Value = "Contend Upload"
$ws.Range("D21").Value = "Content creator user should be able to upload content on site."
$ws.Range("F21").Value = "Login as content creator, attempt by clicking upload content."

# Row 22: new "Subscription Management" test case
$ws.Range("A22").Value = "Subscription_01"
$ws.Range("B22").Value = "Subscription Management"
$ws.Range("D22").Value = "User should be able to subscribe, upgrade or cancel subscription plans."
$ws.Range("F22").Value = "Login as content creator, attempt to change subscription status."

# Row 23: new "Social Media Share" test case
$ws.Range("A23").Value = "Sharing_01"
$ws.Range("B23").Value = "Social Media Share"
$ws.Range("D23").Value = "Share button should be there which lets user share content to social media."
$ws.Range("F23").Value = "Click Share button next to content and attempt to share to social media."

# Row 24 previously held leftover/orphaned text in D24:E24 - clear it entirely so the row disappears
$ws.Range("D24:E24").Clear()

# Update the active selection to reflect where the author left off editing
$ws.Range("H22").Select()
